$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the comment text in C12: the author added a new leading sentence
# about finding red fox data in Peter's file, prefixed onto the previous comment.
$ws.Range("C12").Value = "verkar som att det finns rödrävsdata i Peters fil fram till 2008 (röd text). hur göra detta? Det rimliga vore avstånd till  aktiva rödrävslyor från varje aktiv fjällrävslya. Detta kommer alltså att ändras från år till år. Hur gör jag det i r?"

# Move the active selection from C18 to C12, reflecting where the author was working
$ws.Range("C12").Select()

# Try to reflect the scrolled window position change (best effort; some
# cosmetic window-position metadata may not be persisted by this runtime)
$win = $excel.ActiveWindow
$win.Top = 900
